$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reordered worker rows (rows 16-21) with updated "Valor Mora" (F) / "Salario Basico" (G)
# New ordering: Celso, Jose Luis, Richar, Luis, Roger, Laura
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "73089468"
$ws.Range("D16").Value = "CELSO ANTONIO GONZALEZ FORTICH"
$ws.Range("E16").Value = "2305"
$ws.Range("F16").Value = 156000
$ws.Range("G16").Value = 5000000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73196033"
$ws.Range("D17").Value = "JOSE LUIS CERVANTES MEJIA"
$ws.Range("E17").Value = "1908"
$ws.Range("F17").Value = 1600
$ws.Range("G17").Value = 1200000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047421035"
$ws.Range("D18").Value = "RICHAR RAFAEL SILGADO VILLALOBOS"
$ws.Range("E18").Value = "2305"
$ws.Range("F18").Value = 2667
$ws.Range("G18").Value = 2000000

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143342046"
$ws.Range("D19").Value = "LUIS PUENTES VALLE"
$ws.Range("E19").Value = "2005"
$ws.Range("F19").Value = 68000
$ws.Range("G19").Value = 2000000

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "92226300"
$ws.Range("D20").Value = "ROGER CEDRON RAMIREZ"
$ws.Range("E20").Value = "2312"
$ws.Range("F20").Value = 8000
$ws.Range("G20").Value = 6000000

$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1002392859"
$ws.Range("D21").Value = "LAURA VANESA RODRIGUEZ GONZALEZ"
$ws.Range("E21").Value = "2011"
$ws.Range("F21").Value = 30430
$ws.Range("G21").Value = 0

# Column widths follow the content resize (bestFit) after the data update
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
